# Apply data updates to the "magento_new_users" sheet and refresh the
# selection / window display shown on the first screen (per commit:
# "Modifying display in 1st screen, comments on conftest").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "John"
$ws.Range("C2").Value = "Doe1"
$ws.Range("D2").Value = "Forever24052025@gmail.com"
$ws.Range("E2").Value = "Forever123"
$ws.Range("F2").Value = "Forever123"

# Row 3
$ws.Range("B3").Value = "John"
$ws.Range("C3").Value = "Doe2"
$ws.Range("D3").Value = "Forever26052025@gmail.com"
$ws.Range("E3").Value = "Forever123"
$ws.Range("F3").Value = "Forever123"

# Row 4
$ws.Range("E4").Value = "Forever123"
$ws.Range("F4").Value = "Forever123"

# Row 5
$ws.Range("C5").Value = "Doe3"
$ws.Range("D5").Value = "Forever270420252@gmail.com"
$ws.Range("E5").Value = "Forever123"
$ws.Range("F5").Value = "Forever123"

# Row 6
$ws.Range("B6").Value = "John"
$ws.Range("D6").Value = "Forever270420252@gmail.com"
$ws.Range("E6").Value = "Forever123"
$ws.Range("F6").Value = "Forever123"

# Row 7
$ws.Range("B7").Value = "John"
$ws.Range("C7").Value = "Doe4"
$ws.Range("E7").Value = "Forever123"
$ws.Range("F7").Value = "Forever123"

# Row 8
$ws.Range("E8").Value = "Forev"

# Row 9
$ws.Range("E9").Value = "Forever123"

# Row 10
$ws.Range("B10").Value = "John"
$ws.Range("C10").Value = "Doe5"
$ws.Range("D10").Value = "Forever20042028@gmail.com"
$ws.Range("E10").Value = "Forever123"
$ws.Range("F10").Value = "Forever123*"

# Move/refresh the active cell selection shown on the first screen to G7
[void]$ws.Range("G7").Select()

# Update the workbook window position/size (best-effort; mirrors the
# intended xWindow/yWindow/windowWidth/windowHeight change).
$win = $excel.ActiveWindow
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10300
